$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 105, pushing existing rows 105-123 down to 107-125
$ws.Rows("105:106").Insert()

# Row 105: new Haba "Primera" record
$ws.Cells.Item(105, 1).Value = 9
$ws.Cells.Item(105, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(105, 3).Value = "Metropolitana"
$ws.Cells.Item(105, 4).Value = 44449
$ws.Cells.Item(105, 5).Value = 13
$ws.Cells.Item(105, 6).Value = 100112026
$ws.Cells.Item(105, 7).Value = "Haba"
$ws.Cells.Item(105, 8).Value = "Sin especificar"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 43
$ws.Cells.Item(105, 11).Value = 15000
$ws.Cells.Item(105, 12).Value = 16000
$ws.Cells.Item(105, 13).Value = 15512
$ws.Cells.Item(105, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(105, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(105, 16).Value = 620
$ws.Cells.Item(105, 17).Value = 25
$ws.Cells.Item(105, 18).Value = "Hortaliza"

# Row 106: new Haba "Segunda" record
$ws.Cells.Item(106, 1).Value = 9
$ws.Cells.Item(106, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(106, 3).Value = "Metropolitana"
$ws.Cells.Item(106, 4).Value = 44449
$ws.Cells.Item(106, 5).Value = 13
$ws.Cells.Item(106, 6).Value = 100112026
$ws.Cells.Item(106, 7).Value = "Haba"
$ws.Cells.Item(106, 8).Value = "Sin especificar"
$ws.Cells.Item(106, 9).Value = "Segunda"
$ws.Cells.Item(106, 10).Value = 25
$ws.Cells.Item(106, 11).Value = 13000
$ws.Cells.Item(106, 12).Value = 14000
$ws.Cells.Item(106, 13).Value = 13480
$ws.Cells.Item(106, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(106, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(106, 16).Value = 539
$ws.Cells.Item(106, 17).Value = 25
$ws.Cells.Item(106, 18).Value = "Hortaliza"
